# Updated cryptos list on Mon Sep 18 13:06:35 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row. Numeric-looking Price strings are entered with a leading
# apostrophe (forcing text, matching the original inlineStr cell type) and
# then the cell style is reset to "Normal" so no stray number-format /
# quote-prefix styling is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.434.66"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "1.668.39"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").Value = "'220.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "'19.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.12%  "
$ws.Range("D11").Value = "'0.0850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "1.898.76"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "1.670.32"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "'4.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "'0.536"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'67.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.15%  "
$ws.Range("D17").Value = "27.402.91"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "'224.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.91%  "
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "'6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.93%  "
$ws.Range("D22").Value = "'4.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").Value = "'2.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").Value = "'9.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'147.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").Value = "'7.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.08%  "
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "'16.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("D30").Value = "'0.0516"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").Value = "'1.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").Value = "'3.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "'3.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "1.275.68"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'0.539"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "'0.837"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").Value = "'0.814"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "1.810.51"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").Value = "'2.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.77%  "
$ws.Range("D45").Value = "'62.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").Value = "'92.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").Value = "'1.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").Value = "'0.0518"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "'7.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").Value = "'0.0984"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").Value = "'0.408"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
